$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set cells in an order that reproduces the shared-string table order seen in
# the target workbook: ... , disky mc diskface, Master, DigitalInstantiation.location

# Row 1 - headers (except E1, set later)
$ws.Range("A1").Value = "DigitalInstantiation.filename"
$ws.Range("B1").Value = "Asset.id"
$ws.Range("C1").Value = "DigitalInstantiation.generations"
$ws.Range("D1").Value = "DigitalInstantiation.generations"
$ws.Range("F1").Value = "DigitalInstantition.aapb_preservation_lto"
$ws.Range("G1").Value = "DigitalInstantition.aapb_preservation_disk"

# Row 2 - data
$ws.Range("A2").Value = "sample_digital_instantiation.xml"
$ws.Range("B2").Value = 123456
$ws.Range("C2").Value = "Proxy"
$ws.Range("D2").Value = "Master"
$ws.Range("E2").Value = "American Archive of Public Broadcasting"
$ws.Range("F2").Value = "fhqwhgads"
$ws.Range("G2").Value = "disky mc diskface"

# E1 last, so "DigitalInstantiation.location" is appended to the shared
# strings table after "Master"
$ws.Range("E1").Value = "DigitalInstantiation.location"

# Update selection to match target (E3, as shown in diff's sheetView selection)
$ws.Range("E3").Select()
